# Generate Report for Handoff
# Updates the localization-status workbook to mark d0f1791e-... and
# fe184776-... as "Ready for handoff" with fresh timestamps, and records
# an "Error Detail" note that the handback file version is stale.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyStatus = "Ready for handoff"

$d0f1791eError = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7dec7d23e6033ab8018d708e3ffc521eb9fa6dd1/e2e/d0f1791e-ec4a-48ca-b78a-9dc700635ad5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b02a125ab2f2ff2ad30926a4d8312205ee4427ff/e2e/d0f1791e-ec4a-48ca-b78a-9dc700635ad5.md."
$fe184776Error = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7dec7d23e6033ab8018d708e3ffc521eb9fa6dd1/e2e/fe184776-cbdc-48e7-8954-4414bc65ca3d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b02a125ab2f2ff2ad30926a4d8312205ee4427ff/e2e/fe184776-cbdc-48e7-8954-4414bc65ca3d.md."

# ---------------------------------------------------------------
# Overview sheet: row 4 = d0f1791e-...md, row 5 = fe184776-...md
# ---------------------------------------------------------------
$wsOverview.Range("E4").Value = $readyStatus
$wsOverview.Range("F4").Value = $readyStatus
$wsOverview.Range("G4").Value = "2016-08-25 22:26:18"

$wsOverview.Range("E5").Value = $readyStatus
$wsOverview.Range("F5").Value = $readyStatus
$wsOverview.Range("G5").Value = "2016-08-25 22:26:18"

# ---------------------------------------------------------------
# zh-cn sheet: row 4 = d0f1791e-...md, row 5 = fe184776-...md
# ---------------------------------------------------------------
$wsZhCn.Range("C4").Value = $readyStatus
$wsZhCn.Range("H4").Value = "2016-08-25 22:26:14"
$wsZhCn.Range("P4").Value = $d0f1791eError

$wsZhCn.Range("C5").Value = $readyStatus
$wsZhCn.Range("H5").Value = "2016-08-25 22:26:14"
$wsZhCn.Range("P5").Value = $fe184776Error

$wsZhCn.Columns.Item(16).ColumnWidth = 39.15

# ---------------------------------------------------------------
# de-de sheet: row 4 = d0f1791e-...md, row 5 = fe184776-...md
# ---------------------------------------------------------------
$wsDeDe.Range("C4").Value = $readyStatus
$wsDeDe.Range("H4").Value = "2016-08-25 22:26:18"
$wsDeDe.Range("P4").Value = $d0f1791eError

$wsDeDe.Range("C5").Value = $readyStatus
$wsDeDe.Range("H5").Value = "2016-08-25 22:26:18"
$wsDeDe.Range("P5").Value = $fe184776Error

$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
